$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 12502
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 12502
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H135").Value = 1653.75
$ws.Range("I135").Value = 1675.7142
$ws.Range("K135").Value = 15081.4278
$ws.Range("M135").Value = -12546.4278
$ws.Range("H137").Value = 3183.5
$ws.Range("I137").Value = 2482.1875
$ws.Range("J137").Value = 3985
$ws.Range("K137").Value = 7446.5625
$ws.Range("L137").Value = 11955
$ws.Range("M137").Value = -4896.5625
$ws.Range("N137").Value = -17055
$ws.Range("H138").Value = 4811.2856
$ws.Range("I138").Value = 4315.7
$ws.Range("J138").Value = 5009.52
$ws.Range("K138").Value = 12947.1
$ws.Range("L138").Value = 15028.56
$ws.Range("M138").Value = -7807.099999999999
$ws.Range("N138").Value = -25308.56

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1142.9423
$ws.Range("I32").Value = 1138.66
$ws.Range("K32").Value = 1138.66
$ws.Range("M32").Value = -851.6600000000001
$ws.Range("H53").Value = 30019.5
$ws.Range("J53").Value = 50000
$ws.Range("L53").Value = 50000
$ws.Range("N53").Value = -51364
$ws.Range("H61").Value = 5675.619
$ws.Range("I61").Value = 3895.5293
$ws.Range("J61").Value = 13241
$ws.Range("K61").Value = 3895.5293
$ws.Range("L61").Value = 13241
$ws.Range("M61").Value = -3683.5293
$ws.Range("N61").Value = -13665
$ws.Range("H63").Value = 4484.625
$ws.Range("J63").Value = 7050
$ws.Range("L63").Value = 7050
$ws.Range("N63").Value = -8422
$ws.Range("H66").Value = 4484.625
$ws.Range("J66").Value = 7050
$ws.Range("L66").Value = 35250
$ws.Range("N66").Value = -42114
$ws.Range("H74").Value = 7578797
$ws.Range("I74").Value = 9261668
$ws.Range("K74").Value = 9261668
$ws.Range("M74").Value = -9260794
$ws.Range("H77").Value = 7578797
$ws.Range("I77").Value = 9261668
$ws.Range("K77").Value = 46308340
$ws.Range("M77").Value = -46303972
$ws.Range("H132").Value = 3400.173
$ws.Range("I132").Value = 2120.2708
$ws.Range("K132").Value = 6360.812399999999
$ws.Range("M132").Value = -3830.812399999999
$ws.Range("H136").Value = 5675.619
$ws.Range("I136").Value = 3895.5293
$ws.Range("J136").Value = 13241
$ws.Range("K136").Value = 11686.5879
$ws.Range("L136").Value = 39723
$ws.Range("M136").Value = -9136.5879
$ws.Range("N136").Value = -44823

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 69791.664
$ws.Range("J132").Value = 69791.664
$ws.Range("L132").Value = 69791.664
$ws.Range("N132").Value = -79911.664
$ws.Range("H134").Value = 2100.5
$ws.Range("J134").Value = 7046.3335
$ws.Range("L134").Value = 21139.0005
$ws.Range("N134").Value = -26209.0005

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20581.238
$ws.Range("I31").Value = 2346.074
$ws.Range("J31").Value = 34257.61
$ws.Range("K31").Value = 2346.074
$ws.Range("L31").Value = 34257.61
$ws.Range("M31").Value = -2051.074
$ws.Range("N31").Value = -34847.61
$ws.Range("H34").Value = 20581.238
$ws.Range("I34").Value = 2346.074
$ws.Range("J34").Value = 34257.61
$ws.Range("K34").Value = 2346.074
$ws.Range("L34").Value = 34257.61
$ws.Range("M34").Value = -2144.074
$ws.Range("N34").Value = -34661.61
$ws.Range("H43").Value = 10526
$ws.Range("J43").Value = 10526
$ws.Range("L43").Value = 10526
$ws.Range("N43").Value = -10894
$ws.Range("H58").Value = 4787.095
$ws.Range("I58").Value = 1140.9166
$ws.Range("K58").Value = 1140.9166
$ws.Range("M58").Value = -937.9166
$ws.Range("H62").Value = 8513.0625
$ws.Range("I62").Value = 5071.4
$ws.Range("J62").Value = 14249.167
$ws.Range("K62").Value = 5071.4
$ws.Range("L62").Value = 14249.167
$ws.Range("M62").Value = -4447.4
$ws.Range("N62").Value = -15497.167
$ws.Range("H65").Value = 8513.0625
$ws.Range("I65").Value = 5071.4
$ws.Range("J65").Value = 14249.167
$ws.Range("K65").Value = 25357
$ws.Range("L65").Value = 71245.83499999999
$ws.Range("M65").Value = -22237
$ws.Range("N65").Value = -77485.83499999999
$ws.Range("H88").Value = 18749
$ws.Range("J88").Value = 18749
$ws.Range("L88").Value = 18749
$ws.Range("N88").Value = -19561
$ws.Range("H91").Value = 18749
$ws.Range("J91").Value = 18749
$ws.Range("L91").Value = 18749
$ws.Range("N91").Value = -21557
$ws.Range("H101").Value = 10526
$ws.Range("J101").Value = 10526
$ws.Range("L101").Value = 10526
$ws.Range("N101").Value = -17016
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H132").Value = 2575
$ws.Range("I132").Value = 2035.8158
$ws.Range("K132").Value = 6107.4474
$ws.Range("M132").Value = -3577.4474
$ws.Range("H134").Value = 1613
$ws.Range("I134").Value = 1135.1936
$ws.Range("J134").Value = 5316
$ws.Range("K134").Value = 3405.5808
$ws.Range("L134").Value = 15948
$ws.Range("M134").Value = -870.5808000000002
$ws.Range("N134").Value = -21018
$ws.Range("H136").Value = 4787.095
$ws.Range("I136").Value = 1140.9166
$ws.Range("K136").Value = 3422.7498
$ws.Range("M136").Value = -872.7498000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 5377.778
$ws.Range("H121").Value = 777916.4399999999
$ws.Range("J121").Value = 1072
$ws.Range("L121").Value = 3216
$ws.Range("N121").Value = -5836

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 34000
$ws.Range("J52").Value = 34000
$ws.Range("L52").Value = 34000
$ws.Range("N52").Value = -34518
$ws.Range("H70").Value = 6918
$ws.Range("I70").Value = 6918
$ws.Range("K70").Value = 6918
$ws.Range("M70").Value = -6648
$ws.Range("H73").Value = 6918
$ws.Range("I73").Value = 6918
$ws.Range("K73").Value = 6918
$ws.Range("M73").Value = -5982
$ws.Range("H126").Value = 10157.143
$ws.Range("I126").Value = 1100
$ws.Range("J126").Value = 11666.667
$ws.Range("K126").Value = 3300
$ws.Range("L126").Value = 35000.001
$ws.Range("M126").Value = -830
$ws.Range("N126").Value = -39940.001
$ws.Range("H132").Value = 7551.3
$ws.Range("I132").Value = 4218.5713
$ws.Range("J132").Value = 15327.667
$ws.Range("K132").Value = 12655.7139
$ws.Range("L132").Value = 45983.001
$ws.Range("M132").Value = -10125.7139
$ws.Range("N132").Value = -51043.001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 25395.654
$ws.Range("I22").Value = 23314.578
$ws.Range("J22").Value = 29349.7
$ws.Range("K22").Value = 23314.578
$ws.Range("L22").Value = 29349.7
$ws.Range("M22").Value = -23019.578
$ws.Range("N22").Value = -29939.7
$ws.Range("H27").Value = 25395.654
$ws.Range("I27").Value = 23314.578
$ws.Range("J27").Value = 29349.7
$ws.Range("K27").Value = 23314.578
$ws.Range("L27").Value = 29349.7
$ws.Range("M27").Value = -23207.578
$ws.Range("N27").Value = -29563.7
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H133").Value = 96093.44500000001
$ws.Range("J133").Value = 96093.44500000001
$ws.Range("L133").Value = 96093.44500000001
$ws.Range("N133").Value = -101153.445
$ws.Range("H136").Value = 7450.12
$ws.Range("I136").Value = 3517.3572
$ws.Range("K136").Value = 10552.0716
$ws.Range("M136").Value = -8002.071599999999
$ws.Range("H141").Value = 79998.57000000001
$ws.Range("J141").Value = 79998.57000000001
$ws.Range("L141").Value = 79998.57000000001
$ws.Range("N141").Value = -90358.57000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 46534
$ws.Range("J80").Value = 46534
$ws.Range("L80").Value = 46534
$ws.Range("N80").Value = -48530
$ws.Range("H83").Value = 46534
$ws.Range("J83").Value = 46534
$ws.Range("L83").Value = 139602
$ws.Range("N83").Value = -149586
$ws.Range("H97").Value = 20000
$ws.Range("J97").Value = 20000
$ws.Range("L97").Value = 20000
$ws.Range("N97").Value = -21982
$ws.Range("H132").Value = 4669.2036
$ws.Range("I132").Value = 2307.239
$ws.Range("J132").Value = 18250.5
$ws.Range("K132").Value = 6921.717000000001
$ws.Range("L132").Value = 54751.5
$ws.Range("M132").Value = -4391.717000000001
$ws.Range("N132").Value = -59811.5
$ws.Range("H136").Value = 4081.6775
$ws.Range("I136").Value = 3983.724
$ws.Range("J136").Value = 5502
$ws.Range("K136").Value = 11951.172
$ws.Range("L136").Value = 16506
$ws.Range("M136").Value = -9401.172
$ws.Range("N136").Value = -21606
